$d = $word.ActiveDocument

# 1. Update the day of the month: "1 días del mes de junio del año 2025" -> "3 días del mes de junio del año 2025"
$d.Content.Find.Execute("1 días del mes de junio del año 2025", $true, $false, $false, $false, $false, $true, 1, $false, "3 días del mes de junio del año 2025", 2)

# 2. Update the address: "Eduardo Gomez, 5458" -> "eduardo gomez, calle 3232"
$d.Content.Find.Execute("Eduardo Gomez, 5458", $true, $false, $false, $false, $false, $true, 1, $false, "eduardo gomez, calle 3232", 2)

# 3. Update the RUC number: "4447477-1, " -> "6554878-9, "
$d.Content.Find.Execute("4447477-1, ", $true, $false, $false, $false, $false, $true, 1, $false, "6554878-9, ", 2)

# 4. Update "Campo via" -> "campo via"
$d.Content.Find.Execute(" Campo via", $true, $false, $false, $false, $false, $true, 1, $false, " campo via", 2)
